# "app is stable - pulling from excel and adding to page"
#
# Rebuilds the Date / Link table:
#   - header row: A1 "Date", B1 "Link" (was "Email Links")
#   - row 2: 3/18/2024 + the PrimaryArms email-link URL (now a real hyperlink)
#   - row 3 (new): 3/19/2024 + the OpticsPlanet email-link URL (new hyperlink)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url1 = "https://cors-anywhere.herokuapp.com/https://news.primaryarms.com/q/R2wAIrCyKnByoGxG1BC9dPcdY_eYGio1HHvBFO66scis4S-QZq4QmAXD2"
$url2 = "https://cors-anywhere.herokuapp.com/https://www.opticsplanet.com/view/email/791098eaeb6c3a8981d1917c1dbc3bb0"

# Wipe the old rows outright (not just ClearFormats) so the stale border/
# wrap/Arial-font formatting and the ht/thickBot row flags don't linger -
# the new table uses plain default formatting plus a real date format and
# the built-in Hyperlink style.
$ws.Rows("1:2").Delete()

$ws.Range("A1").Value = "Date"

$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Value = 45369
$ws.Range("B2").Value = $url1
$ws.Hyperlinks.Add($ws.Range("B2"), $url1)

$ws.Range("A3").NumberFormat = "mm-dd-yy"
$ws.Range("A3").Value = 45370
$ws.Range("B3").Value = $url2
$ws.Hyperlinks.Add($ws.Range("B3"), $url2)

$ws.Range("B1").Value = "Link"

$ws.Columns("A").ColumnWidth = 8.592447916666666
$ws.Columns("B").ColumnWidth = 117.73697916666667

[void]$ws.Range("B11").Select()
